$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.644.16'
$ws.Range("E2").Value = '  -1.44%  '

$ws.Range("D3").Value = '2.239.52'
$ws.Range("E3").Value = '  -1.56%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '115.38'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +3.49%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '285.41'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +7.89%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.628'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -3.52%  '

$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.614'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.97%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '46.75'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +0.38%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0931'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '9.17'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("E13").Value = '  -2.71%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '15.41'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +1.05%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.883'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +2.96%  '

$ws.Range("D16").Value = '2.576.12'
$ws.Range("E16").Value = '  -1.56%  '

$ws.Range("D17").Value = '2.244.49'
$ws.Range("E17").Value = '  -1.73%  '

$ws.Range("D18").Value = '42.999.47'
$ws.Range("E18").Value = '  -0.29%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.0000108'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '6.90'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +2.76%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '73.01'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +1.28%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '3.20'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +11.57%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '2.36'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '231.89'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '9.19'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -1.71%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '12.09'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +6.48%  '

$ws.Range("E27").Value = '  -1.65%  '

$ws.Range("E28").Value = '  -0.81%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '40.30'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '2.24'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("B31").Value = 'WEMIXToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '3.30'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -1.51%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '175.60'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +1.57%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '21.19'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -1.12%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.0902'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +0.71%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '4.64'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +19.72%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '5.59'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.80%  '

$ws.Range("E37").Value = '  -2.82%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.0373'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -1.40%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '4.63'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -0.66%  '

$ws.Range("E40").Value = '  +1.61%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.62'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +1.89%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '72.71'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -2.35%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '13.50'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -5.21%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.234'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("E45").Value = '  +0.33%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '1.34'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '5.60'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -7.82%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '1.30'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +2.79%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '8.55'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.656'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +9.27%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.473'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +9.37%  '
